$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace every occurrence of "Both" with "SM" across the used range
# (these are the Class=Global rows that used to be tagged "Both").
$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Value2 -eq "Both") {
        $cell.Value2 = "SM"
    }
}

# Update the active selection to reflect the latest edit position
$ws.Range("H9").Select()
